# fix IJPP and add two journals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 130: Journal of Communication ---
$ws.Range("E130").Value = "https://academic.oup.com/joc/pages/General_Instructions#Submission%20Guidelines"
$ws.Range("A130").Value = "JOURNAL OF COMMUNICATION"
$ws.Range("B130").Value = "35 pages"
$ws.Range("F130").Value = 20190918

$ws.Hyperlinks.Add($ws.Range("E130"), "https://academic.oup.com/joc/pages/General_Instructions", "Submission%20Guidelines")
$ws.Range("E130").Style = "Link"

# --- Row 131: Communication Research ---
$ws.Range("E131").Value = "https://us.sagepub.com/en-us/nam/journal/communication-research#submission-guidelines"
$ws.Range("A131").Value = "COMMUNICATION RESEARCH"
$ws.Range("B131").Value = "35 pages"
$ws.Range("F131").Value = 20190918

$ws.Hyperlinks.Add($ws.Range("E131"), "https://us.sagepub.com/en-us/nam/journal/communication-research", "submission-guidelines")
$ws.Range("E131").Style = "Link"

# Keep the active selection in sync with where Excel would leave the cursor
$null = $ws.Range("A131").Select()
